$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.587.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.358.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.46%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'513.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.30%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'127.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.29%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.373.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -5.83%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0958"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.95%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -8.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.316"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.85%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.782.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'56.527.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.79%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -4.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.305.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -8.53%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'309.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -5.52%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.91%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'64.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.389"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.99%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.466.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.73%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.63%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.13%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'173.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -5.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0₃0714"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.41%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.42%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E38").Value = "'  -6.44%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -7.46%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.803"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.31%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'35.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.88%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'4.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.59%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Mantle"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.569"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.23%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'122.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'251.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -10.47%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0905"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0487"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.75%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0207"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.83%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'16.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.47%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Applied changes"